$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right after "2021-Q4" (before "总计"),
#    holding the per-fund detail for that quarter.
# ---------------------------------------------------------------------------
$ws2021Q4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $ws2021Q4)
$newSheet.Name = "2022-Q1"

# Reuse the header/index-column styling already used on the other quarter
# sheets (bold+border header row, bold+border A column).
$ws2021Q4.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$ws2021Q4.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

$newSheet.Cells.Item(2, 1).Value = 0

# B2:G2 hold text-looking numbers (fund code / sizes / ratios) that must be
# stored as text, not auto-converted to numeric cells.
$textRange = $newSheet.Range("B2:G2")
$textRange.NumberFormat = "@"
$newSheet.Cells.Item(2, 2).Value = "004209"
$newSheet.Cells.Item(2, 3).Value = "大成智惠量化多策略灵活配置混合"
$newSheet.Cells.Item(2, 4).Value = "1.26"
$newSheet.Cells.Item(2, 5).Value = "94.24"
$newSheet.Cells.Item(2, 6).Value = "5.94"
$newSheet.Cells.Item(2, 7).Value = "0.0748"
$textRange.ClearFormats()

$newSheet.Cells.Item(2, 8).Value = 8

# ---------------------------------------------------------------------------
# 2) Add the corresponding "2022-Q1" row to the "总计" summary sheet.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows(2).Insert()

# The inherited row format isn't used by data rows on this sheet, only the
# A column carries the (bold+border) index style.
$wsTotal.Range("A2:D2").ClearFormats()
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 1
$wsTotal.Cells.Item(2, 4).Value = 0.07000000000000001

# Renumber the running index (column A) of the rows that got pushed down.
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(4, 1).Value = 2
$wsTotal.Cells.Item(5, 1).Value = 3
$wsTotal.Cells.Item(6, 1).Value = 4
$wsTotal.Cells.Item(7, 1).Value = 5
